# Update column F (dSF) values for specific rows to reflect re-pulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F29").Value = 0
$ws.Range("F37").Value = 2
$ws.Range("F41").Value = -2
$ws.Range("F59").Value = 2
$ws.Range("F60").Value = 3
$ws.Range("F74").Value = -5
